$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Test: Copy B5:J19 to L5 using Copy(Destination)
$ws.Range("B5:J19").Copy($ws.Range("L5"))
